$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A2").Value = 100083269
$ws.Range("B2").Value = 77506
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 569103.6928627139
$ws.Range("R2").Value = 6903955.500130533

$ws.Range("A3").Value = 100083992
$ws.Range("B3").Value = 92501
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 1004672
$ws.Range("F3").Value = "Källmossor"
$ws.Range("G3").Value = "Philonotis"
$ws.Range("H3").Value = "Brid."
$ws.Range("Q3").Value = 569096.4477078969
$ws.Range("R3").Value = 6903967.465335313

$ws.Range("A4").Value = 100083268
$ws.Range("B4").Value = 78569
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("Q4").Value = 569295.646553766
$ws.Range("R4").Value = 6904060.092043239

$ws.Range("A5").Value = 100083273
$ws.Range("B5").Value = 77506
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 568877.6041835379
$ws.Range("R5").Value = 6904149.805104633

$ws.Range("A6").Value = 100083266
$ws.Range("B6").Value = 77506
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 569294.2761642485
$ws.Range("R6").Value = 6904058.666035361

$ws.Range("A7").Value = 100083272
$ws.Range("B7").Value = 77506
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 569076.6973299341
$ws.Range("R7").Value = 6904043.004845407

$ws.Range("A8").Value = 100083991
$ws.Range("B8").Value = 92688
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 2387
$ws.Range("F8").Value = "Källpraktmossa"
$ws.Range("G8").Value = "Pseudobryum cinclidioides"
$ws.Range("H8").Value = "(Huebener) T.J.Kop."
$ws.Range("Q8").Value = 569211.1671905057
$ws.Range("R8").Value = 6904083.976534305
